$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: add "Before" / "After" labels
$ws.Range("D7").Value2 = "Before"
$ws.Range("H7").Value2 = "After"

# Row 8: add "Audit and batch changes" label (inherits row's customFormat style)
$ws.Range("H8").Value2 = "Audit and batch changes"

# Row 7: add "Large volume" label (written after H8 so shared-string order matches)
$ws.Range("L7").Value2 = "Large volume"

# Row 9: start markers + elapsed time value
$ws.Range("G9").Value2 = "start"
$ws.Range("K9").Value2 = "start"
$ws.Range("L9").NumberFormat = "h:mm:ss"
$ws.Range("L9").Value2 = 0.3704513888888889

# Row 20: start marker + elapsed time value
$ws.Range("G20").Value2 = "start"
$ws.Range("L20").NumberFormat = "h:mm:ss"
$ws.Range("L20").Value2 = 0.37048611111111113

# Row 21: end marker + elapsed time value
$ws.Range("G21").Value2 = "end"
$ws.Range("L21").NumberFormat = "h:mm:ss"
$ws.Range("L21").Value2 = 0.37248842592592596

# Row 22: end marker + elapsed time value
$ws.Range("G22").Value2 = "end"
$ws.Range("L22").NumberFormat = "h:mm:ss"
$ws.Range("L22").Value2 = 0.37248842592592596

# Row 24: before/after summary values
$ws.Range("L24").NumberFormat = "h:mm:ss"
$ws.Range("L24").Value2 = "2 min"
$ws.Range("M24").Value2 = "57 secs"

# Row 26: audit/batch summary values
$ws.Range("K26").Value2 = "value"
$ws.Range("L26").Value2 = "2 min"
$ws.Range("M26").Value2 = "5 sec"

# Leave selection on K26, matching the final saved state
$ws.Range("K26").Select() | Out-Null
